$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.601.30'
$ws.Range("E2").Value = '  +2.11%  '
$ws.Range("D3").Value = '1.689.26'
$ws.Range("E3").Value = '  +3.39%  '
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '221.06'
$ws.Range("E5").Value = '  +2.77%  '
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("D8").Value = '31.12'
$ws.Range("E8").Value = '  +4.58%  '
$ws.Range("D9").Value = '0.267'
$ws.Range("E9").Value = '  +2.51%  '
$ws.Range("D10").Value = '0.0628'
$ws.Range("E10").Value = '  +2.40%  '
$ws.Range("D11").Value = '0.0904'
$ws.Range("E11").Value = '  -1.38%  '
$ws.Range("D12").Value = '1.931.48'
$ws.Range("E12").Value = '  +3.42%  '
$ws.Range("D13").Value = '10.83'
$ws.Range("E13").Value = '  +13.61%  '
$ws.Range("D14").Value = '0.623'
$ws.Range("E14").Value = '  +8.52%  '
$ws.Range("D15").Value = '1.686.01'
$ws.Range("E15").Value = '  +3.27%  '
$ws.Range("D16").Value = '4.03'
$ws.Range("E16").Value = '  +3.46%  '
$ws.Range("D17").Value = '30.598.11'
$ws.Range("E17").Value = '  +2.00%  '
$ws.Range("D18").Value = '66.09'
$ws.Range("E18").Value = '  +1.78%  '
$ws.Range("D19").Value = '248.21'
$ws.Range("E19").Value = '  -0.34%  '
$ws.Range("D20").Value = '0.0₃0721'
$ws.Range("E20").Value = '  +1.87%  '
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("E22").Value = '  +3.49%  '
$ws.Range("D23").Value = '10.22'
$ws.Range("E23").Value = '  +5.63%  '
$ws.Range("D24").Value = '2.18'
$ws.Range("E24").Value = '  +2.51%  '
$ws.Range("D25").Value = '157.47'
$ws.Range("E25").Value = '  -1.56%  '
$ws.Range("E26").Value = '  +1.67%  '
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("D28").Value = '6.79'
$ws.Range("E28").Value = '  +2.25%  '
$ws.Range("E30").Value = '  +2.14%  '
$ws.Range("E31").Value = '  +1.18%  '
$ws.Range("E32").Value = '  +3.57%  '
$ws.Range("E33").Value = '  +2.92%  '
$ws.Range("D34").Value = '1.507.72'
$ws.Range("E34").Value = '  +5.13%  '
$ws.Range("E35").Value = '  +5.70%  '
$ws.Range("E36").Value = '  -0.75%  '
$ws.Range("E37").Value = '  +4.63%  '
$ws.Range("D38").Value = '79.63'
$ws.Range("E38").Value = '  +8.23%  '
$ws.Range("E39").Value = '  +5.03%  '
$ws.Range("D40").Value = '2.70'
$ws.Range("E40").Value = '  -5.56%  '
$ws.Range("E41").Value = '  +1.44%  '
$ws.Range("D42").Value = '0.856'
$ws.Range("E42").Value = '  +2.72%  '
$ws.Range("D43").Value = '2.02'
$ws.Range("E43").Value = '  +1.62%  '
$ws.Range("D44").Value = '0.0503'
$ws.Range("E44").Value = '  +1.24%  '
$ws.Range("E45").Value = '  -1.98%  '
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("D47").Value = '52.32'
$ws.Range("E47").Value = '  -5.01%  '
$ws.Range("D48").Value = '1.824.09'
$ws.Range("E48").Value = '  +2.75%  '
$ws.Range("E49").Value = '  -0.73%  '
$ws.Range("D50").Value = '95.61'
$ws.Range("D51").Value = '0.0₆0115'
$ws.Range("E51").Value = '  +5.86%  '
